# Applies the "Updated cryptos list" data refresh to Sheet1 (rows 2-51).
# Numeric-looking Price strings (single decimal point, e.g. "303.51") would
# be auto-coerced to floating point numbers by a plain Value assignment, so
# those specific cells are forced to Text format first, then the format bump
# is cleared again so the cell keeps its original (unstyled) appearance.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.909.12"
$ws.Range("E2").Value = "  -1.37%  "
$ws.Range("D3").Value = "2.338.06"
$ws.Range("E3").Value = "  +0.01%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "303.51"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "94.20"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -4.29%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.501"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -1.27%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.494"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -2.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.08"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -4.80%  "
$ws.Range("E11").Value = "  -2.26%  "
$ws.Range("E12").Value = "  -4.51%  "
$ws.Range("E13").Value = "  +1.94%  "
$ws.Range("E14").Value = "  -3.21%  "
$ws.Range("D15").Value = "2.703.67"
$ws.Range("E15").Value = "  +0.50%  "
$ws.Range("D16").Value = "2.334.28"
$ws.Range("E16").Value = "  -1.63%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.793"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.18%  "
$ws.Range("D18").Value = "42.851.37"
$ws.Range("E18").Value = "  -1.27%  "
$ws.Range("E19").Value = "  -5.57%  "
$ws.Range("E20").Value = "  +2.03%  "
$ws.Range("D21").Value = "0.0₃0887"
$ws.Range("E21").Value = "  -1.69%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.91"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.22%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.78"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.92%  "
$ws.Range("E24").Value = "  -1.66%  "
$ws.Range("E25").Value = "  -0.07%  "
$ws.Range("E26").Value = "  -1.67%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.54"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.36%  "
$ws.Range("E28").Value = "  -0.30%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.11"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.40%  "
$ws.Range("E30").Value = "  -6.66%  "
$ws.Range("E31").Value = "  +0.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0763"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +8.35%  "
$ws.Range("E33").Value = "  -1.18%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "17.23"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -4.24%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.37"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -3.58%  "
$ws.Range("E36").Value = "  -1.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "125.84"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -23.54%  "
$ws.Range("E38").Value = "  +1.93%  "
$ws.Range("E39").Value = "  -0.45%  "
$ws.Range("E40").Value = "  -1.63%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "22.10"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +21.13%  "
$ws.Range("E42").Value = "  -1.52%  "
$ws.Range("D43").Value = "1.933.68"
$ws.Range("E43").Value = "  -3.03%  "
$ws.Range("E44").Value = "  -0.47%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.17"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -5.69%  "
$ws.Range("E46").Value = "  +1.10%  "
$ws.Range("E47").Value = "  -3.46%  "
$ws.Range("B48").Value = "HuobiToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.88"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.58%  "
$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value = "2.567.13"
$ws.Range("E49").Value = "  +0.41%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "52.68"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.83%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "71.63"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.79%  "
